$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.988.40"
$ws.Range("E2").Value = "  +0.97%  "
$ws.Range("D3").Value = "2.790.25"
$ws.Range("E3").Value = "  -0.91%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'358.92"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.93%  "
$ws.Range("D6").Value = "'109.73"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.80%  "
$ws.Range("D7").Value = "'0.567"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.59%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").Value = "'0.596"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.48%  "
$ws.Range("D10").Value = "'40.23"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.64%  "
$ws.Range("D11").Value = "'0.0855"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.48%  "
$ws.Range("E12").Value = "  +1.05%  "
$ws.Range("D13").Value = "'19.56"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.95%  "
$ws.Range("D14").Value = "'7.62"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.83%  "
$ws.Range("D15").Value = "3.225.00"
$ws.Range("E15").Value = "  -1.31%  "
$ws.Range("D16").Value = "2.788.35"
$ws.Range("D17").Value = "'0.935"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.20%  "
$ws.Range("D18").Value = "51.895.55"
$ws.Range("E18").Value = "  +0.99%  "
$ws.Range("E19").Value = "  +1.07%  "
$ws.Range("D20").Value = "'3.14"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.77%  "
$ws.Range("D21").Value = "'13.10"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.91%  "
$ws.Range("D22").Value = "0.0₃0981"
$ws.Range("E22").Value = "  -0.84%  "
$ws.Range("D23").Value = "'274.51"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.59%  "
$ws.Range("D24").Value = "'70.29"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.09%  "
$ws.Range("D25").Value = "'2.74"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.40%  "
$ws.Range("D26").Value = "'26.72"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.47%  "
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("D28").Value = "'10.20"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.67%  "
$ws.Range("E29").Value = "  +5.24%  "
$ws.Range("E30").Value = "  -1.21%  "
$ws.Range("D31").Value = "'0.0466"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.93%  "
$ws.Range("D32").Value = "'51.54"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.90%  "
$ws.Range("D33").Value = "'34.15"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.66%  "
$ws.Range("D34").Value = "'5.73"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.20%  "
$ws.Range("D35").Value = "'0.0845"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.05%  "
$ws.Range("D36").Value = "'5.28"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.60%  "
$ws.Range("E37").Value = "  -0.07%  "
$ws.Range("E38").Value = "  +1.93%  "
$ws.Range("D39").Value = "'18.10"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.65%  "
$ws.Range("E40").Value = "  -2.72%  "
$ws.Range("D41").Value = "'2.59"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.67%  "
$ws.Range("E42").Value = "  -1.35%  "
$ws.Range("E43").Value = "  -3.09%  "
$ws.Range("D44").Value = "'2.25"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.13%  "
$ws.Range("D45").Value = "'21.97"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -7.32%  "
$ws.Range("D46").Value = "2.078.13"
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("D47").Value = "'3.26"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.38%  "
$ws.Range("E48").Value = "  -5.58%  "
$ws.Range("D49").Value = "'5.73"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.59%  "
$ws.Range("D50").Value = "'0.943"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.07%  "
$ws.Range("D51").Value = "'8.98"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.90%  "
